# Apply the updated "split" dataset to Sheet1.
# Rows 2-12 get new values (shifted window of data), and rows 13-17 are
# newly appended rows that continue the same table / style as column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for columns A..K (HOURS, param_P_to_demand1, param_E_pv1_solar,
# param_E_solar_th1_solar, param_E_pvt1_solar, param_Q_to_demand1,
# param_net1_cost_sell_electric, param_net1_cost_buy_electric,
# param_net1_cost_sell_thermal, param_net1_cost_buy_thermal)
# for worksheet rows 2 through 17.
$data = @(
    @(19, 20, 433.5207499755035, 0, 0, 0, 650.2811249632553, 0.38, 0.41, 0.17, 0.33),
    @(20, 21, 393.8044005019047, 0, 0, 0, 590.7066007528571, 0.35, 0.49, 0.16, 0.37),
    @(21, 22, 333.3011855326704, 0, 0, 0, 499.9517782990056, 0.3, 0.4, 0.19, 0.34),
    @(22, 23, 282.2034344188766, 0, 0, 0, 423.3051516283149, 0.39, 0.41, 0.14, 0.3),
    @(23, 24, 208.1255268683671, 0, 0, 0, 312.1882903025506, 0.32, 0.45, 0.18, 0.39),
    @(24, 25, 165.1488800499209, 0, 0, 0, 247.7233200748814, 0.4, 0.48, 0.2, 0.4),
    @(25, 26, 130.2933375268219, 0, 0, 0, 195.4400062902328, 0.37, 0.46, 0.1, 0.33),
    @(26, 27, 103.6774555661557, 0, 0, 0, 155.5161833492336, 0.38, 0.48, 0.19, 0.32),
    @(27, 28, 97.78916096758275, 0, 0, 0, 146.6837414513741, 0.35, 0.43, 0.14, 0.32),
    @(28, 29, 94.78573553475324, 0, 0, 0, 142.1786033021299, 0.4, 0.49, 0.14, 0.39),
    @(29, 30, 98.71785172641819, 0, 0, 0, 148.0767775896273, 0.37, 0.46, 0.15, 0.33),
    @(30, 31, 130.7873219730109, 0, 0, 0, 196.1809829595164, 0.37, 0.49, 0.13, 0.4),
    @(31, 32, 201.8420447128433, 0, 0, 0, 302.7630670692649, 0.37, 0.46, 0.19, 0.31),
    @(32, 33, 294.8296768634683, 0.00277778, 0.00277778, 0.00555556, 442.2445152952025, 0.34, 0.46, 0.18, 0.38),
    @(33, 34, 349.0296502993302, 0.0138889, 0.0138889, 0.0277778, 523.5444754489954, 0.39, 0.5, 0.15, 0.4),
    @(34, 35, 364.1653337305623, 0.03055558, 0.03055558, 0.06111116, 546.2480005958435, 0.38, 0.43, 0.11, 0.37)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $colNum = $j + 1
        $ws.Cells.Item($rowNum, $colNum).Value = $rowVals[$j]
    }
}

# Column A (the HOURS index column) uses the bordered/bold/centered style
# already present on A2:A12. Propagate that same style onto the newly
# created A13:A17 cells so the whole column stays visually consistent.
$ws.Range("A2").Copy()
$ws.Range("A13:A17").PasteSpecial(-4122)
$excel.CutCopyMode = 0
